$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "23.961.60"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value2 = "  -1.91%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "1.648.30"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value2 = "  -1.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value2 = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "310.18"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value2 = "  -1.07%  "

$ws.Range("E6").Value2 = "  +0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.3886"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value2 = "  -2.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.3805"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value2 = "  -2.84%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "51.97"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value2 = "  -0.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "1.344"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value2 = "  -4.02%  "

$ws.Range("E11").Value2 = "  +0.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.08432"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value2 = "  -2.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "23.89"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value2 = "  -2.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "7.050"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value2 = "  -3.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "8.044"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value2 = "  +2.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "0.00001306"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value2 = "  -4.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "1.650.70"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value2 = "  -0.78%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "94.09"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value2 = "  -1.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "0.06988"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value2 = "  -0.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "19.64"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value2 = "  -4.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "6.942"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value2 = "  -0.94%  "

$ws.Range("E22").Value2 = "  +0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "13.73"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value2 = "  -0.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "23.948.08"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value2 = "  -1.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "2.460"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value2 = "  +1.41%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "2.949"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value2 = "  -3.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "22.03"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value2 = "  -2.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "153.42"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value2 = "  -2.66%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "5.391"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value2 = "  -0.97%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "138.07"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value2 = "  -3.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "7.812"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value2 = "  -3.64%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "2.502"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value2 = "  -0.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "1.830.97"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value2 = "  -0.79%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "1.013"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value2 = "  -5.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "0.08109"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value2 = "  -1.99%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "6.746"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value2 = "  -2.36%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.02930"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value2 = "  -3.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.2676"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value2 = "  -3.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "10.70"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value2 = "  -3.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.09075"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value2 = "  -1.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.7578"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value2 = "  -2.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "13.38"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value2 = "  -3.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "1.419"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value2 = "  -1.75%  "

$ws.Range("E44").Value2 = "  -2.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.6941"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value2 = "  -2.62%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "2.447"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value2 = "  -3.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "4.091"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value2 = "  -1.10%  "

$ws.Range("E48").Value2 = "  +0.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.08299"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value2 = "  -1.81%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "134.08"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value2 = "  -1.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "1.230"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value2 = "  -3.41%  "
